$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 6: path separator question/answer.
# B6 uses the wrapped-text style like the other question cells (B2/B4).
$ws.Range("B6").Value = "路径分隔"
$ws.Range("B6").WrapText = $true

# New row 11 (answer cell) must be added to the shared-string table before
# the "File.separator" text so the resulting shared string order matches.
$ws.Range("C11").Value = "a"

# C6 holds the corresponding answer text.
$ws.Range("C6").Value = "File.separator"

# Update the active sheet's view/selection to match the new scroll position
# and the new active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("C7").Select()
$ws.Range("C7").Activate()
